$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-22 Tuesday", "2024-10-23 Wednesday"),
    @("399×3=", "268×7="),
    @("198×5=", "571×4="),
    @("716×5=", "659×2="),
    @("332×5=", "285×6="),
    @("981×5=", "575×8="),
    @("551×2=", "841×3="),
    @("802×8=", "677×6="),
    @("588×3=", "628×6="),
    @("977×8=", "974×4="),
    @("215×3=", "190×8="),
    @("380×5=", "405×4="),
    @("520×7=", "991×4="),
    @("384×8=", "476×9="),
    @("125×9=", "229×4="),
    @("200×8=", "837×7="),
    @("675×7=", "440×6="),
    @("825×2=", "749×3="),
    @("600×6=", "131×4="),
    @("818×3=", "264×7="),
    @("633×6=", "254×9="),
    @("863×9=", "981×4="),
    @("783×5=", "438×9="),
    @("293×4=", "277×3="),
    @("311×3=", "685×2="),
    @("301×7=", "996×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
